$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "64.069.83"
Set-TextValue "E2" "  -0.91%  "
Set-TextValue "D3" "3.405.27"
Set-TextValue "E3" "  -1.15%  "
Set-TextValue "E4" "  -0.01%  "
Set-TextValue "D5" "573.13"
Set-TextValue "E5" "  +0.04%  "
Set-TextValue "D6" "162.67"
Set-TextValue "E6" "  +2.24%  "
Set-TextValue "E7" "  +0.02%  "
Set-TextValue "D8" "3.403.36"
Set-TextValue "E8" "  -1.15%  "
Set-TextValue "E9" "  -5.39%  "
Set-TextValue "E10" "  +1.24%  "
Set-TextValue "E11" "  -1.97%  "
Set-TextValue "D13" "3.992.13"
Set-TextValue "E13" "  -1.16%  "
Set-TextValue "E14" "  +0.36%  "
Set-TextValue "D15" "26.82"
Set-TextValue "E15" "  -2.53%  "
Set-TextValue "E16" "  -0.47%  "
Set-TextValue "D17" "64.089.66"
Set-TextValue "E17" "  -0.89%  "
Set-TextValue "D18" "3.408.04"
Set-TextValue "E18" "  -1.40%  "
Set-TextValue "E19" "  -1.09%  "
Set-TextValue "E20" "  -2.71%  "
Set-TextValue "D21" "373.52"
Set-TextValue "E21" "  -1.87%  "
Set-TextValue "E22" "  -2.59%  "
Set-TextValue "E23" "  +0.01%  "
Set-TextValue "D24" "70.40"
Set-TextValue "E24" "  -2.80%  "
Set-TextValue "E25" "  -3.67%  "
Set-TextValue "E26" "  -3.30%  "
Set-TextValue "D27" "9.47"
Set-TextValue "E27" "  -4.49%  "
Set-TextValue "E28" "  -0.80%  "
Set-TextValue "D29" "1.00"
Set-TextValue "E29" "  +1.20%  "
Set-TextValue "E30" "  -0.57%  "
Set-TextValue "E31" "  -3.07%  "
Set-TextValue "E32" "  -1.02%  "
Set-TextValue "D33" "0.999"
Set-TextValue "E33" "  +0.05%  "
Set-TextValue "D34" "22.78"
Set-TextValue "E34" "  -1.98%  "
Set-TextValue "D35" "7.02"
Set-TextValue "E35" "  -0.06%  "
Set-TextValue "E36" "  -6.55%  "
Set-TextValue "D37" "159.06"
Set-TextValue "E37" "  -1.36%  "
Set-TextValue "E38" "  +7.20%  "
Set-TextValue "E39" "  -3.06%  "
Set-TextValue "E40" "  -3.00%  "
Set-TextValue "D41" "25.76"
Set-TextValue "E41" "  -2.41%  "
Set-TextValue "D42" "42.62"
Set-TextValue "E42" "  -0.78%  "
Set-TextValue "D43" "2.723.86"
Set-TextValue "E43" "  -5.43%  "
Set-TextValue "E44" "  -1.36%  "
Set-TextValue "D45" "25.74"
Set-TextValue "E45" "  -1.01%  "
Set-TextValue "D47" "0.0303"
Set-TextValue "E47" "  -2.12%  "
Set-TextValue "E48" "  -1.44%  "
Set-TextValue "D49" "328.51"
Set-TextValue "E50" "  -2.57%  "
Set-TextValue "E51" "  -2.06%  "
